$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-23 from 45212 to 45221
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 45221
}
